# hello world.xlsx — CKP-R edit
# * rename the worksheet tab to "CKP-R"
# * narrow column C (Uraian Kegiatan) from 80 to 65 characters
# * add a new "Pegawai yang Dinilai / Pejabat Penilai" signature block
#   (rows 19-25) with centered, wrap-text, merged B:C / H:I cells
# * move the selection onto the new block

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename sheet -----------------------------------------------------
$ws.Name = "CKP-R"

# --- narrow the "Uraian Kegiatan" column -------------------------------
# ColumnWidth is expressed in "characters of the Normal style font" and the
# engine round-trips it into the stored OOXML width with a fixed +5/6
# offset, so back the request off by that amount to land exactly on 65.
$ws.Range("C1").ColumnWidth = 65 - (5/6)

# --- new signature block (rows 19-25) ----------------------------------
$block = $ws.Range("B19:I25")
$block.HorizontalAlignment = -4108   # xlCenter
$block.VerticalAlignment = -4108     # xlCenter
$block.WrapText = $true

$ws.Range("B19").Value = "Pegawai yang Dinilai"
$ws.Range("H19").Value = "Pejabat Penilai"
$ws.Range("B24").Value = "Indra"
$ws.Range("H24").Value = "Iva"

$ws.Range("B19:C19").Merge()
$ws.Range("H19:I19").Merge()
$ws.Range("B20:C20").Merge()
$ws.Range("H20:I20").Merge()
$ws.Range("B21:C21").Merge()
$ws.Range("H21:I21").Merge()
$ws.Range("B22:C22").Merge()
$ws.Range("H22:I22").Merge()
$ws.Range("B23:C23").Merge()
$ws.Range("H23:I23").Merge()
$ws.Range("B24:C24").Merge()
$ws.Range("H24:I24").Merge()
$ws.Range("B25:C25").Merge()
$ws.Range("H25:I25").Merge()

# --- move the selection onto the new block -----------------------------
$null = $ws.Range("B19:I25").Select()
